# Update file with new data
# Adds a new data row (row 4) for IP 85.104.3.240 with its VirusTotal
# link/hyperlink, analysis stats, country and AS owner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 values -----------------------------------------------------
$ws.Range("C4").Value = "Malicious"
$ws.Range("D4").Value = "https://www.virustotal.com/gui/ip-address/85.104.3.240/detection"
$ws.Range("E4").Value = "{'harmless': 56, 'malicious': 12, 'suspicious': 1, 'undetected': 21, 'timeout': 0}"
$ws.Range("F4").Value = "Turkey"
$ws.Range("G4").Value = 45306.97109953704
$ws.Range("H4").Value = 45340.0480787037
$ws.Range("I4").Value = "Turk Telekom"

# Match the date number formatting used by the other rows.
$ws.Range("G4").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("H4").NumberFormat = $ws.Range("H2").NumberFormat

# Hyperlink for the Link column, same as rows 2 and 3.
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.virustotal.com/gui/ip-address/85.104.3.240/detection") | Out-Null

# Restore the "Link" column style (hyperlink look) after Hyperlinks.Add.
$ws.Range("D4").Style = $ws.Range("D2").Style
